$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 8 (pushes existing rows 8..37 down to 9..38),
# carrying the row-8 date style down with the shift.
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the latest weekly price record.
$ws.Cells.Item(8, 1).Value = 11
$ws.Cells.Item(8, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(8, 3).Value = "Bíobío"
$ws.Cells.Item(8, 4).Value = 44462
$ws.Cells.Item(8, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(8, 5).Value = 8
$ws.Cells.Item(8, 6).Value = 100112001
$ws.Cells.Item(8, 7).Value = "Berenjena"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 100
$ws.Cells.Item(8, 11).Value = 7000
$ws.Cells.Item(8, 12).Value = 7500
$ws.Cells.Item(8, 13).Value = 7250
$ws.Cells.Item(8, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(8, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(8, 16).Value = 121
$ws.Cells.Item(8, 17).Value = 60
$ws.Cells.Item(8, 18).Value = "Hortaliza"
